$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 42.31746133333333
$ws.Range("N2").Value = 126.952384
$ws.Range("O2").Value = 0.6904142182914543
$ws.Range("P2").Value = 0.6904142182914543
$ws.Range("Q2").Value = 4.903973112433778
$ws.Range("R2").Value = 44.135758011904
$ws.Range("S2").Value = 0.6904142182914543
$ws.Range("T2").Value = 0.6904142182914543

# Row 3
$ws.Range("M3").Value = 11.08476666666667
$ws.Range("N3").Value = 33.2543
$ws.Range("O3").Value = 0.1808492350906109
$ws.Range("P3").Value = 0.1808492350906109
$ws.Range("Q3").Value = 1.284561880088889
$ws.Range("R3").Value = 11.5610569208
$ws.Range("S3").Value = 0.1808492350906109
$ws.Range("T3").Value = 0.1808492350906109

# Row 4
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.140061
$ws.Range("N4").Value = 0.420183
$ws.Range("O4").Value = 0.002285111223152439
$ws.Range("P4").Value = 0.002285111223152439
$ws.Range("Q4").Value = 0.016231015672
$ws.Range("R4").Value = 0.146079141048
$ws.Range("S4").Value = 0.002285111223152439
$ws.Range("T4").Value = 0.002285111223152439

# Row 5
$ws.Range("M5").Value = 7.750570000000001
$ws.Range("N5").Value = 23.25171
$ws.Range("O5").Value = 0.1264514353947823
$ws.Range("P5").Value = 0.1264514353947823
$ws.Range("Q5").Value = 0.8981773879733335
$ws.Range("R5").Value = 8.083596491760002
$ws.Range("S5").Value = 0.1264514353947823
$ws.Range("T5").Value = 0.1264514353947823
